$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers (row 1) are unchanged.

# Data rows: columns A-D are categorical (Sending cluster, Ligand symbol,
# Receptor symbol, Target cluster); columns E-T are numeric metrics.
$rows = @(
    @{ Row=2;  A="FAPs"; B="Wnt1"; C="Ryk"; D="FAPs";
       E=1; F=0.3333333333333333; G=0.01948966666666667; H=0.058469;
       I=0.07096062449330311; J=0.07096062449330311; K=3; L=1;
       M=9.182053333333334; N=27.54616; O=0.1632474062119586; P=0.1844076677824049;
       Q=0.1789551587822222; R=1.61059642904; S=0.01158413789171251; T=0.01308568326719302 },
    @{ Row=3;  A="FAPs"; B="Wnt1"; C="Ryk"; D="ECs";
       E=1; F=0.3333333333333333; G=0.01948966666666667; H=0.058469;
       I=0.07096062449330311; J=0.07096062449330311; K=3; L=1;
       M=27.25159233333333; N=81.75477699999999; O=0.4845051103561108; P=0.5473070568326256;
       Q=0.5311244507125554; R=4.780120056413; S=0.03438078520106636; T=0.03883725054243484 },
    @{ Row=4;  A="FAPs"; B="Wnt1"; C="Ryk"; D="M1";
       E=1; F=0.3333333333333333; G=0.01948966666666667; H=0.058469;
       I=0.07096062449330311; J=0.07096062449330311; K=3; L=1;
       M=0.3190983333333333; N=0.957295; O=0.00567323814751954; P=0.006408607890528383;
       Q=0.006219120150555555; R=0.055972081355; S=0.0004025765218472166; T=0.0004547588180446039 },
    @{ Row=5;  A="FAPs"; B="Wnt1"; C="Ryk"; D="M2";
       E=1; F=0.3333333333333333; G=0.01948966666666667; H=0.058469;
       I=0.07096062449330311; J=0.07096062449330311; K=2; L=0.6666666666666666;
       M=0.1312043333333333; N=0.393613; O=0.00233267726976492; P=0.002635040794754541;
       Q=0.002557128721888889; R=0.023014158497; S=0.000165528235803852; T=0.000186984140361112 },
    @{ Row=6;  A="FAPs"; B="Wnt1"; C="Ryk"; D="sCs";
       E=1; F=0.3333333333333333; G=0.01948966666666667; H=0.058469;
       I=0.07096062449330311; J=0.07096062449330311; K=2; L=1;
       M=19.362295; N=38.72459; O=0.344241568014646; P=0.2592416266996866;
       Q=0.3773646754516666; R=2.26418805271; S=0.02442759664287316; T=0.01839594772526952 },
    @{ Row=7;  A="ECs";  B="Wnt1"; C="Ryk"; D="FAPs";
       E=3; F=1; G=0.255165; H=0.7654949999999999;
       I=0.9290393755066968; J=0.9290393755066969; K=3; L=1;
       M=9.182053333333334; N=27.54616; O=0.1632474062119586; P=0.1844076677824049;
       Q=2.3429386388; R=21.0864477492; S=0.151663268320246; T=0.1713219845152119 },
    @{ Row=8;  A="ECs";  B="Wnt1"; C="Ryk"; D="ECs";
       E=3; F=1; G=0.255165; H=0.7654949999999999;
       I=0.9290393755066968; J=0.9290393755066969; K=3; L=1;
       M=27.25159233333333; N=81.75477699999999; O=0.4845051103561108; P=0.5473070568326256;
       Q=6.953652557734999; R=62.58287301961499; S=0.4501243251550444; T=0.5084698062901908 },
    @{ Row=9;  A="ECs";  B="Wnt1"; C="Ryk"; D="M1";
       E=3; F=1; G=0.255165; H=0.7654949999999999;
       I=0.9290393755066968; J=0.9290393755066969; K=3; L=1;
       M=0.3190983333333333; N=0.957295; O=0.00567323814751954; P=0.006408607890528383;
       Q=0.08142272622499999; R=0.732804536025; S=0.005270661625672322; T=0.005953849072483779 },
    @{ Row=10; A="ECs";  B="Wnt1"; C="Ryk"; D="M2";
       E=3; F=1; G=0.255165; H=0.7654949999999999;
       I=0.9290393755066968; J=0.9290393755066969; K=2; L=0.6666666666666666;
       M=0.1312043333333333; N=0.393613; O=0.00233267726976492; P=0.002635040794754541;
       Q=0.033478753715; R=0.3013087834349999; S=0.002167149033961068; T=0.002448056654393429 },
    @{ Row=11; A="ECs";  B="Wnt1"; C="Ryk"; D="sCs";
       E=3; F=1; G=0.255165; H=0.7654949999999999;
       I=0.9290393755066968; J=0.9290393755066969; K=2; L=1;
       M=19.362295; N=38.72459; O=0.344241568014646; P=0.2592416266996866;
       Q=4.940580003675; R=29.64348002205; S=0.3198139713717728; T=0.2408456789744171 }
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
}
